$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iAU_TC_ID_33"
$ws.Range("B2").Value = "@RegressionA Validation of Create Question (MCQ/Type A)-Negative Scenario"
$ws.Range("C2").Value = "passed"
